$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 8

$ws.Cells.Item($row, 1).Value = 9

# Column B holds a date-like string ("2026-02-16"); Excel's normal entry
# path would auto-convert this to a date serial. Force text entry, then
# drop back to the workbook's default ("Normal") cell style so no stray
# number-format style is left behind on the cell.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "21:22:06"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"
$ws.Cells.Item($row, 6).Value = 69348.815
# Exit Price (G) stays blank/open, same as the other OPEN trades above it.
$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.7199
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.072% move"
# Exit Reason (M) stays blank, same as the other OPEN trades above it.
$ws.Cells.Item($row, 14).Value = 0
